# Zeitaufzeichnung.xlsx - Protokoll erweitert
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stedronsky")

# New protocol rows (5-7)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Design Vorlage"
$ws.Range("C5").Value = 42042

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Desgin Umsetzung und RM "
$ws.Range("C6").Value = 42044

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Fehler überarbeitung FK"
$ws.Range("C7").Value = 42045

$ws.Range("C5:C7").NumberFormat = "m/d/yyyy"

# D3: mark the "DatabaseMetaData" task as done by "Stedronsky"
$ws.Range("D3").Value = "Stedronsky"

# leave the active selection where the author left off
$ws.Range("H16").Select()
